$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45177 -> 45178, i.e. 2023-09-08 -> 2023-09-09) for every data row,
# rows 2 through 420.
for ($r = 2; $r -le 420; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}
